# Reorder the "Recorded By" (column G) values so that "System" (exact case)
# appears first in the comma-separated list, while preserving the relative
# order of the remaining entries.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7th column
    $val = $cell.Value2

    if ($null -eq $val) { continue }
    if ($val -isnot [string]) { continue }
    if ($val.IndexOf(",") -lt 0) { continue }

    $parts = $val -split ", "
    $trimmed = @()
    foreach ($p in $parts) { $trimmed += $p.Trim() }

    $hasSystem = $false
    foreach ($p in $trimmed) {
        if ($p.Equals("System")) { $hasSystem = $true }
    }

    if ($hasSystem) {
        $rest = @()
        foreach ($p in $trimmed) {
            if (-not $p.Equals("System")) { $rest += $p }
        }
        $newParts = @("System") + $rest
        $newVal = [string]::Join(", ", $newParts)

        if (-not $newVal.Equals($val)) {
            $cell.Value = $newVal
        }
    }
}
